$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N held sample-specific text values like 60/30", 100/5", etc.
# Replace them with plain numeric values (drops the now-unused shared strings).
$ws.Range("N11").Value = 20
$ws.Range("N12").Value = 20
$ws.Range("N13").Value = 25
$ws.Range("N14").Value = 50
$ws.Range("N17").Value = 10
$ws.Range("N18").Value = 20
$ws.Range("N19").Value = 17
$ws.Range("N20").Value = 100
$ws.Range("N21").Value = 50

# Update the sheet's last active selection, as recorded by Excel on save.
$ws.Range("R24").Select()
